$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the "github for win" shared string into B5 (previously reverted to 2222)
$ws.Range("B5").Value = "github for win"

# Move the selection on to B6, matching the post-edit cursor position
$ws.Range("B6").Select()
